# "add the original paper" - Robot Schedule.xlsx
# - Remove the old D-column "TOTAL" header/formula (D1/D2)
# - Correct several HOURS values in column C
# - Add a new "Totoal" grand-total row (row 17) with =SUM(C2:C16)
# - Leave selection on F12, as in the authored workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old TOTAL header cell and its per-row SUM formula.
$ws.Range("D1").ClearContents()
$ws.Range("D2").ClearContents()

# Corrected HOURS figures.
$ws.Range("C3").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 6
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 5
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 2

# New grand-total row.
$ws.Range("A17").Value = "Totoal"
$ws.Range("C17").Formula = "=SUM(C2:C16)"

# Restore the author's final selection.
$ws.Range("F12").Select() | Out-Null
